# Generate Report for Handback
# The f0820fe7-175a-4264-a194-9d9298b01209.md file has now been handed back
# in sync with en-US (for both zh-cn and de-de locales). Update the status,
# the handback datetime, and clear the stale error detail accordingly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f0820fe7-175a-4264-a194-9d9298b01209.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for f0820fe7-175a-4264-a194-9d9298b01209.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-30 16:57:36"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row for f0820fe7-175a-4264-a194-9d9298b01209.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-30 16:57:44"
$wsDeDe.Range("P3").Value = ""
